$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Step 1: append " html" (spell-checked as two extra runs) to the
#             last row's "Mejora prototipo ciudad" paragraph ---
$lastRow = $t.Rows($t.Rows.Count)
$para1 = $lastRow.Cells(2).Range.Paragraphs(1).Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mejora prototipo ciudad</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>html</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para1.InsertXML($xml1)

# --- Step 2: add a new table row for 19/08/2020 with three bullet items ---
$newRow = $t.Rows.Add()
$c1 = $newRow.Cells(1)
$c2 = $newRow.Cells(2)

$c1.Range.Text = "19/08/2020"

$c2para = $c2.Range.Paragraphs(1).Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Acabado prototipo ciudad</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Mejora c&#243;digo prototipo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y ciudad </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>html</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Retoque </w:t></w:r><w:r><w:t xml:space="preserve">prototipo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$c2para.InsertXML($xml2)
